$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update cell C10 on the sheet: value changed from 18 to 1
$ws.Range("C10").Value = 1
